# OPEN APNEE LYON - POST-EVENT
# Update the "Battement Epreuve" (G column) figures on the Données sheet
# for the sta / dwf / dnf events (rows 2, 3 and 5). These feed the
# Visuel_Planning schedule sheet, whose TIME()-based formulas recompute
# automatically from the new values.
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Données")
$wsData.Range("G2").Value = 17
$wsData.Range("G3").Value = 29
$wsData.Range("G5").Value = 14

# Visuel_Planning: update the active selection and flip the print
# orientation to landscape for the post-event layout.
$wsPlanning = $wb.Worksheets.Item("Visuel_Planning")
$wsPlanning.Activate()
$wsPlanning.Range("H2").Select() | Out-Null
$wsPlanning.PageSetup.Orientation = 2
